$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.283.20"
$ws.Range("E2").Value = "  -1.60%  "
$ws.Range("D3").Value = "1.827.64"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("E4").Value = "  -0.76%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.83%  "
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4248"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.95%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3716"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07256"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8649"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.12"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.40%  "
$ws.Range("D12").Value = "1.832.29"
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.733"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.324"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.82%  "
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008881"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.004"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("B20").Value = "BitDAO"
$ws.Range("C20").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.5044"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.69%  "
$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "27.359.25"
$ws.Range("E22").Value = "  -1.37%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.137"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.41%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.63%  "
$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "2.053.34"
$ws.Range("E25").Value = "  -1.57%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.992"
$ws.Range("D26").Style = "Normal"
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.177"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.77%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.249"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.04%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "116.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.48%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08863"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.16%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.200"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.89%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7586"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.59%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.464"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.47%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.807"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.95%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.005"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.122"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01980"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05277"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.355"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.59%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.872"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1699"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.86%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5067"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.37%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.698"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.72%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "107.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.81%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4761"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.005"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.80%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06398"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.71%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.675"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.41%  "
